$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 10.08677976748346
$ws.Range("F2").Value = 1.982163051966514
$ws.Range("G2").Value = 0.1936338431263439
$ws.Range("E3").Value = 9.635115734468156
$ws.Range("F3").Value = 1.876121013669253
$ws.Range("G3").Value = 0.2297411636881241
$ws.Range("E4").Value = 9.002853604636213
$ws.Range("F4").Value = 1.7175899044724
$ws.Range("G4").Value = 0.2802860150204473
$ws.Range("E5").Value = 10.789857869266
$ws.Range("F5").Value = 2.507317193046481
$ws.Range("G5").Value = 0.1374277595213425
$ws.Range("E6").Value = 10.626151874059
$ws.Range("F6").Value = 2.4895996542852
$ws.Range("G6").Value = 0.1505148871532741
$ws.Range("E7").Value = 9.762707854168175
$ws.Range("F7").Value = 2.259172142577704
$ws.Range("G7").Value = 0.2195410830299205
$ws.Range("E8").Value = 11.40201970920793
$ws.Range("F8").Value = 2.640696912204143
$ws.Range("G8").Value = 0.08848978311682376
$ws.Range("E9").Value = 11.44121351274578
$ws.Range("F9").Value = 2.638265343091528
$ws.Range("G9").Value = 0.08535651784677545
$ws.Range("E10").Value = 18.30459480980873
$ws.Range("F10").Value = 3.712976353278274
$ws.Range("G10").Value = -0.4633219035372529
$ws.Range("E11").Value = 10.01257175343638
$ws.Range("F11").Value = 2.054612337237419
$ws.Range("G11").Value = 0.1995662449905412
$ws.Range("E12").Value = 9.638744627291498
$ws.Range("F12").Value = 1.978685536274874
$ws.Range("G12").Value = 0.2294510595689584
$ws.Range("E13").Value = 8.924712824351674
$ws.Range("F13").Value = 1.801773200800383
$ws.Range("G13").Value = 0.2865328135176521
$ws.Range("E14").Value = 10.9643123099862
$ws.Range("F14").Value = 2.563426274932812
$ws.Range("G14").Value = 0.1234813702716649
$ws.Range("E15").Value = 11.34103483544833
$ws.Range("F15").Value = 2.692996430858489
$ws.Range("G15").Value = 0.09336508915250041
$ws.Range("E16").Value = 11.02053569852597
$ws.Range("F16").Value = 2.587557177463951
$ws.Range("G16").Value = 0.1189867110455982
$ws.Range("E17").Value = 11.40201970920793
$ws.Range("F17").Value = 2.640696912204143
$ws.Range("G17").Value = 0.08848978311682376
$ws.Range("E18").Value = 11.44121351274578
$ws.Range("F18").Value = 2.638265343091528
$ws.Range("G18").Value = 0.08535651784677545
$ws.Range("E19").Value = 23.13974408159993
$ws.Range("F19").Value = 4.349686086511766
$ws.Range("G19").Value = -0.8498576291187243
$ws.Range("E20").Value = 10.27456678128292
$ws.Range("F20").Value = 2.26940482829159
$ws.Range("G20").Value = 0.1786216096763381
$ws.Range("E21").Value = 9.871645667994313
$ws.Range("F21").Value = 2.175874192038457
$ws.Range("G21").Value = 0.2108322811825379
$ws.Range("E22").Value = 9.046923328262753
$ws.Range("F22").Value = 1.957795048889328
$ws.Range("G22").Value = 0.2767629546885686
$ws.Range("E23").Value = 10.9643123099862
$ws.Range("F23").Value = 2.563426274932812
$ws.Range("G23").Value = 0.1234813702716649
$ws.Range("E24").Value = 11.34103483544833
$ws.Range("F24").Value = 2.692996430858489
$ws.Range("G24").Value = 0.09336508915250041
$ws.Range("E25").Value = 13.14682448941866
$ws.Range("F25").Value = 2.98812208732648
$ws.Range("G25").Value = -0.05099492434639141
$ws.Range("E26").Value = 11.40201970920793
$ws.Range("F26").Value = 2.640696912204143
$ws.Range("G26").Value = 0.08848978311682376
$ws.Range("E27").Value = 11.44121351274578
$ws.Range("F27").Value = 2.638265343091528
$ws.Range("G27").Value = 0.08535651784677545
$ws.Range("E28").Value = 27.6574986234277
$ws.Range("F28").Value = 4.912789929284439
$ws.Range("G28").Value = -1.211019908019262
